# Automatic update of files.
# Update the "Förändrad" date column (C) for all data rows from 46076 to 46077,
# and re-sync the re-ordered records in rows 7-10 and 12-14 (A/B/G columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column C: bump every data row's "Förändrad" value by one day (46076 -> 46077) ---
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = 46077
}

# --- Rows whose record (Beteckning/Datum/Area) moved to a different row ---

# Row 7 <- old row 10
$ws.Cells.Item(7, 1).Value = "A 19922-2025"
$ws.Cells.Item(7, 2).Value = 45771.63034722222
$ws.Cells.Item(7, 7).Value = 10.1

# Row 8 <- old row 9
$ws.Cells.Item(8, 1).Value = "A 25015-2023"
$ws.Cells.Item(8, 2).Value = 45085.6989699074
$ws.Cells.Item(8, 7).Value = 1.8

# Row 9 <- old row 7
$ws.Cells.Item(9, 1).Value = "A 62884-2021"
$ws.Cells.Item(9, 2).Value = 44504
$ws.Cells.Item(9, 7).Value = 0.8

# Row 10 <- old row 13
$ws.Cells.Item(10, 1).Value = "A 14271-2021"
$ws.Cells.Item(10, 2).Value = 44278
$ws.Cells.Item(10, 7).Value = 6.7

# Row 12 <- old row 14
$ws.Cells.Item(12, 1).Value = "A 25634-2025"
$ws.Cells.Item(12, 2).Value = 45803.59570601852
$ws.Cells.Item(12, 7).Value = 6

# Row 13 <- old row 12
$ws.Cells.Item(13, 1).Value = "A 3402-2026"
$ws.Cells.Item(13, 2).Value = 46042.39047453704
$ws.Cells.Item(13, 7).Value = 5.5

# Row 14 <- old row 8
$ws.Cells.Item(14, 1).Value = "A 28266-2025"
$ws.Cells.Item(14, 2).Value = 45818.56381944445
$ws.Cells.Item(14, 7).Value = 1.9
